# rnaSample_hbrown_08.08.19.xlsx -- "major accuracy check update"
#
# 1. Fix the RNA-extraction-reagent typo "Trizol" -> "TRIzol" for every
#    sample row (the G column, rows 2-27 all share this value).
# 2. Give the whole G column the same font G2 already used (Arial 11,
#    black) so every sample row is visually consistent.
# 3. Make the C1 header cell ("bioSampleNumber") match the font already
#    used by the rest of the header row (Arial 10, black).
# 4. Widen column H so the "roboticRNAPrep" header is not cramped.
# 5. Tighten the data-row height now that the column formatting is more
#    consistent.
# 6. Leave the active cell on H2, where the reviewer finished checking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- 1. Typo fix, applied across the whole sample column at once --------
$ws.Range("G2:G27").Replace("Trizol", "TRIzol")

# -- 2. Normalize the sample column's font to match the first row -------
$ws.Range("G2:G27").Font.Name = "Arial"
$ws.Range("G2:G27").Font.Size = 11
$ws.Range("G2:G27").Font.Color = 0

# -- 3. Normalize the C1 header cell's font to match the rest of row 1 --
$ws.Range("C1").Font.Name = "Arial"
$ws.Range("C1").Font.Size = 10
$ws.Range("C1").Font.Color = 0

# -- 4. Widen column H ---------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 17.14

# -- 5. Shrink the data rows' height back down ---------------------------
$ws.Range("A2:Z27").RowHeight = 15

# -- 6. Leave the active selection on H2 ---------------------------------
$ws.Range("H2").Select()
